$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "Done"
$ws.Range("C3").Value = "Halfway done"

$ws.Range("C3").Select()
